$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (P3 header) - SuppliersPartNumber was blank, now filled in
$ws.Range("E3").Value = "C880557"

# Row 4 (P5 header) - BOM line item updated to a new part
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = "Header,Male Pin,Covered(4 Sided) 2 1 0.138""（3.50mm） 2 P=3.5mm Pluggable System Terminal Block RoHS"
$ws.Range("D4").Value = "SL_3_50_180G_02"
$ws.Range("E4").Value = "C192777"
